# Add team record (Wins/Losses/Ties) columns to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns AD/AE/AF ---
# Copy the style of the last existing header cell (AC1, style index 1 = bold/centered/bordered)
# onto the new header cells so they match the rest of the header row.
$ws.Range("AC1").Copy($ws.Range("AD1:AF1"))

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-47): team record repeated on every row ---
$wins = 106
$losses = 56
$ties = 0

for ($r = 2; $r -le 47; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # column AD = 30
    $ws.Cells.Item($r, 31).Value = $losses  # column AE = 31
    $ws.Cells.Item($r, 32).Value = $ties    # column AF = 32
}
